{"js": "/* Positional replacement of all 101 run texts (1 heading paragraph +\n   100 table-cell paragraphs, in document order) per the diff. Some old\n   values repeat (e.g. \"4+40=\" appears twice with different replacements),\n   so matching must be done by position, not by Find/Replace on text. */\n\nconst oldValues = [\n  \"2023-07-08 Saturday\",\n  \"29+70=\",\n  \"3+62=\",\n  \"25+21=\",\n  \"40+30=\",\n  \"79-1=\",\n  \"88-67=\",\n  \"61-31=\",\n  \"86-39=\",\n  \"26+9=\",\n  \"7+66=\",\n  \"43+3=\",\n  \"6+70=\",\n  \"83-12=\",\n  \"66+16=\",\n  \"80-1=\",\n  \"14-12=\",\n  \"86-77=\",\n  \"98-79=\",\n  \"21-16=\",\n  \"62+27=\",\n  \"86-46=\",\n  \"31-14=\",\n  \"65+13=\",\n  \"41-32=\",\n  \"22-9=\",\n  \"44+52=\",\n  \"49-17=\",\n  \"83-71=\",\n  \"63+28=\",\n  \"7+64=\",\n  \"25+64=\",\n  \"71-62=\",\n  \"79-63=\",\n  \"44+37=\",\n  \"4+40=\",\n  \"20-3=\",\n  \"23+63=\",\n  \"0+38=\",\n  \"84-32=\",\n  \"74-14=\",\n  \"72-61=\",\n  \"2+13=\",\n  \"46+52=\",\n  \"25-1=\",\n  \"13+49=\",\n  \"23-17=\",\n  \"37+27=\",\n  \"12-10=\",\n  \"56-50=\",\n  \"41-25=\",\n  \"70+13=\",\n  \"12+20=\",\n  \"74-9=\",\n  \"19+30=\",\n  \"65-15=\",\n  \"88-7=\",\n  \"73-60=\",\n  \"1+88=\",\n  \"86+11=\",\n  \"67+18=\",\n  \"42+40=\",\n  \"69+15=\",\n  \"61-58=\",\n  \"86-75=\",\n  \"71+21=\",\n  \"39+43=\",\n  \"35+2=\",\n  \"96-37=\",\n  \"81-67=\",\n  \"4+40=\",\n  \"56-42=\",\n  \"72+18=\",\n  \"50+14=\",\n  \"75-70=\",\n  \"27+34=\",\n  \"51+12=\",\n  \"20+40=\",\n  \"52+32=\",\n  \"72+12=\",\n  \"76-57=\",\n  \"42-29=\",\n  \"57-20=\",\n  \"78-68=\",\n  \"70-7=\",\n  \"30-14=\",\n  \"62-60=\",\n  \"87-16=\",\n  \"67+22=\",\n  \"66+8=\",\n  \"8+75=\",\n  \"63-34=\",\n  \"48-4=\",\n  \"7+89=\",\n  \"2+49=\",\n  \"94-57=\",\n  \"60-6=\",\n  \"52+36=\",\n  \"72-60=\",\n  \"55-19=\",\n  \"40+48=\"\n];\nconst newValues = [\n  \"2023-07-09 Sunday\",\n  \"9+13=\",\n  \"41+20=\",\n  \"45+9=\",\n  \"14+79=\",\n  \"30+23=\",\n  \"38-37=\",\n  \"88-16=\",\n  \"32-14=\",\n  \"23+24=\",\n  \"56+36=\",\n  \"84+14=\",\n  \"64+11=\",\n  \"26+72=\",\n  \"69-15=\",\n  \"80-78=\",\n  \"45+49=\",\n  \"8+58=\",\n  \"0+11=\",\n  \"28-0=\",\n  \"10+11=\",\n  \"95-86=\",\n  \"53-7=\",\n  \"47-28=\",\n  \"98-10=\",\n  \"63-29=\",\n  \"46+6=\",\n  \"52-4=\",\n  \"23+49=\",\n  \"79-10=\",\n  \"97-1=\",\n  \"2+3=\",\n  \"93-41=\",\n  \"17+6=\",\n  \"41+46=\",\n  \"46-38=\",\n  \"45+34=\",\n  \"85-41=\",\n  \"26+61=\",\n  \"55-25=\",\n  \"91-75=\",\n  \"12+70=\",\n  \"14+81=\",\n  \"94-33=\",\n  \"51-44=\",\n  \"44+39=\",\n  \"68+16=\",\n  \"57-4=\",\n  \"11+49=\",\n  \"57-36=\",\n  \"24+47=\",\n  \"88+6=\",\n  \"40-16=\",\n  \"5+60=\",\n  \"3+14=\",\n  \"44+28=\",\n  \"13+61=\",\n  \"8+2=\",\n  \"57+0=\",\n  \"71-39=\",\n  \"64-62=\",\n  \"56+0=\",\n  \"56+29=\",\n  \"21+37=\",\n  \"37-25=\",\n  \"95-32=\",\n  \"70-66=\",\n  \"52-13=\",\n  \"45+39=\",\n  \"82-62=\",\n  \"62-23=\",\n  \"21+31=\",\n  \"91-48=\",\n  \"29+56=\",\n  \"87-6=\",\n  \"45-8=\",\n  \"21+48=\",\n  \"61+11=\",\n  \"47+46=\",\n  \"77-72=\",\n  \"9+1=\",\n  \"66+1=\",\n  \"2+10=\",\n  \"7+62=\",\n  \"25+41=\",\n  \"49+40=\",\n  \"16+78=\",\n  \"85-16=\",\n  \"56+20=\",\n  \"60-55=\",\n  \"13+16=\",\n  \"65-37=\",\n  \"16+18=\",\n  \"83-22=\",\n  \"95-65=\",\n  \"48-13=\",\n  \"93-88=\",\n  \"80-10=\",\n  \"24+58=\",\n  \"64-8=\",\n  \"90-40=\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== newValues.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + newValues.length +\n    \" got \" + paragraphs.items.length\n  );\n}\n\n// Load text for every paragraph so we can sanity-check before writing.\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text !== oldValues[i]) {\n    throw new Error(\n      \"Paragraph \" + i + \" text mismatch: expected \" + JSON.stringify(oldValues[i]) +\n      \" got \" + JSON.stringify(p.text)\n    );\n  }\n  if (newValues[i] !== oldValues[i]) {\n    p.insertText(newValues[i], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Positional replacement of all 101 run texts (1 heading paragraph +\n# 100 table-cell paragraphs, in document order) per the diff. Some old\n# values repeat (e.g. \"4+40=\" appears twice with different replacements),\n# so matching must be done by position, not by Find/Replace on text.\n#\n# $d.Paragraphs also yields an extra empty \"end of row\" marker paragraph\n# after the last cell of every table row (Range.Text = \"`r`a\", i.e. a\n# carriage return followed by the cell-mark char 7) - those are skipped\n# below so the remaining paragraphs line up 1:1 with the diff order.\n\n$oldValues = @(\n  '2023-07-08 Saturday',\n  '29+70=',\n  '3+62=',\n  '25+21=',\n  '40+30=',\n  '79-1=',\n  '88-67=',\n  '61-31=',\n  '86-39=',\n  '26+9=',\n  '7+66=',\n  '43+3=',\n  '6+70=',\n  '83-12=',\n  '66+16=',\n  '80-1=',\n  '14-12=',\n  '86-77=',\n  '98-79=',\n  '21-16=',\n  '62+27=',\n  '86-46=',\n  '31-14=',\n  '65+13=',\n  '41-32=',\n  '22-9=',\n  '44+52=',\n  '49-17=',\n  '83-71=',\n  '63+28=',\n  '7+64=',\n  '25+64=',\n  '71-62=',\n  '79-63=',\n  '44+37=',\n  '4+40=',\n  '20-3=',\n  '23+63=',\n  '0+38=',\n  '84-32=',\n  '74-14=',\n  '72-61=',\n  '2+13=',\n  '46+52=',\n  '25-1=',\n  '13+49=',\n  '23-17=',\n  '37+27=',\n  '12-10=',\n  '56-50=',\n  '41-25=',\n  '70+13=',\n  '12+20=',\n  '74-9=',\n  '19+30=',\n  '65-15=',\n  '88-7=',\n  '73-60=',\n  '1+88=',\n  '86+11=',\n  '67+18=',\n  '42+40=',\n  '69+15=',\n  '61-58=',\n  '86-75=',\n  '71+21=',\n  '39+43=',\n  '35+2=',\n  '96-37=',\n  '81-67=',\n  '4+40=',\n  '56-42=',\n  '72+18=',\n  '50+14=',\n  '75-70=',\n  '27+34=',\n  '51+12=',\n  '20+40=',\n  '52+32=',\n  '72+12=',\n  '76-57=',\n  '42-29=',\n  '57-20=',\n  '78-68=',\n  '70-7=',\n  '30-14=',\n  '62-60=',\n  '87-16=',\n  '67+22=',\n  '66+8=',\n  '8+75=',\n  '63-34=',\n  '48-4=',\n  '7+89=',\n  '2+49=',\n  '94-57=',\n  '60-6=',\n  '52+36=',\n  '72-60=',\n  '55-19=',\n  '40+48='\n)\n$newValues = @(\n  '2023-07-09 Sunday',\n  '9+13=',\n  '41+20=',\n  '45+9=',\n  '14+79=',\n  '30+23=',\n  '38-37=',\n  '88-16=',\n  '32-14=',\n  '23+24=',\n  '56+36=',\n  '84+14=',\n  '64+11=',\n  '26+72=',\n  '69-15=',\n  '80-78=',\n  '45+49=',\n  '8+58=',\n  '0+11=',\n  '28-0=',\n  '10+11=',\n  '95-86=',\n  '53-7=',\n  '47-28=',\n  '98-10=',\n  '63-29=',\n  '46+6=',\n  '52-4=',\n  '23+49=',\n  '79-10=',\n  '97-1=',\n  '2+3=',\n  '93-41=',\n  '17+6=',\n  '41+46=',\n  '46-38=',\n  '45+34=',\n  '85-41=',\n  '26+61=',\n  '55-25=',\n  '91-75=',\n  '12+70=',\n  '14+81=',\n  '94-33=',\n  '51-44=',\n  '44+39=',\n  '68+16=',\n  '57-4=',\n  '11+49=',\n  '57-36=',\n  '24+47=',\n  '88+6=',\n  '40-16=',\n  '5+60=',\n  '3+14=',\n  '44+28=',\n  '13+61=',\n  '8+2=',\n  '57+0=',\n  '71-39=',\n  '64-62=',\n  '56+0=',\n  '56+29=',\n  '21+37=',\n  '37-25=',\n  '95-32=',\n  '70-66=',\n  '52-13=',\n  '45+39=',\n  '82-62=',\n  '62-23=',\n  '21+31=',\n  '91-48=',\n  '29+56=',\n  '87-6=',\n  '45-8=',\n  '21+48=',\n  '61+11=',\n  '47+46=',\n  '77-72=',\n  '9+1=',\n  '66+1=',\n  '2+10=',\n  '7+62=',\n  '25+41=',\n  '49+40=',\n  '16+78=',\n  '85-16=',\n  '56+20=',\n  '60-55=',\n  '13+16=',\n  '65-37=',\n  '16+18=',\n  '83-22=',\n  '95-65=',\n  '48-13=',\n  '93-88=',\n  '80-10=',\n  '24+58=',\n  '64-8=',\n  '90-40='\n)\n\n$d = $word.ActiveDocument\n$total = $d.Paragraphs.Count\n\n$idx = 0\nfor ($i = 1; $i -le $total; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $txt = $r.Text\n\n    # Skip table row-end mark paragraphs (\"`r`a\").\n    if ($txt.Length -eq 2 -and [int][char]$txt[1] -eq 7) {\n        continue\n    }\n\n    # Range excluding the trailing paragraph mark.\n    $content = $d.Range($r.Start, $r.End - 1)\n\n    if ($content.Text -ne $oldValues[$idx]) {\n        throw \"Paragraph at position $idx text mismatch: expected '$($oldValues[$idx])' got '$($content.Text)'\"\n    }\n\n    if ($newValues[$idx] -ne $oldValues[$idx]) {\n        $content.Text = $newValues[$idx]\n    }\n\n    $idx = $idx + 1\n}\n\nif ($idx -ne $newValues.Length) {\n    throw \"Unexpected paragraph count: expected $($newValues.Length) got $idx\"\n}\n"}
